$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "Search values through autocomplete tab" "Search values through autocomplete tab*"
Replace-Text "Sort by (name, place, category, sleeps, price)" "Sort by (name, place, category, sleeps, price)*"
Replace-Text "See properties" "Show properties"
Replace-Text "View bookings" "View bookings*"
Replace-Text "Sort by (name)" "Sort by (name)*"
Replace-Text "See guests" "Show guests"
Replace-Text "Sort by (property name, guest name, check in/out date)" "Sort by (property name, booking platform, guest name, check in/out date)*"
Replace-Text "See bookings" "Show bookings"
